$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe.
$ws.Range("D2").Value = "26.932.30"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.875.28"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("D14").Value = "1.860.87"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "26.969.95"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "2.079.40"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E39").Value = "  +5.49%  "
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("E51").Value = "  +0.70%  "

# Numeric-looking values that must remain plain text (matches original inlineStr cells):
# Temporarily format as Text, assign, then restore style so no stray numFmt/style lingers.
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "306.42"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9992"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5156"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3728"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07192"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.8990"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "20.72"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07558"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "94.87"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "5.257"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.9990"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.000008514"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "14.27"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.9991"
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "5.033"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.40"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.426"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "146.16"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "18.02"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.114"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "114.79"
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.765"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.09178"
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.05030"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.7529"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.992"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.172"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.260"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01994"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.5578"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.491"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.072"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "6.578"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.717"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "115.79"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.4774"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.9990"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.564"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "37.09"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "63.51"
$cell.Style = "Normal"
